$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reactions")

# Add new column I data. Set I2 first then I1 so that the shared-string
# table receives "fdh" before "Enzyme ID" (matching the target order).
$ws.Range("I2").Value = "fdh"
$ws.Range("I1").Value = "Enzyme ID"

# Give I1 ("Enzyme ID" header) the same look as the other header cells
# (bold font, centered/top aligned) by copying A1's format, then trim the
# border down to just the left/right thin edges.
$ws.Range("A1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Borders(7).LineStyle = 1
$ws.Range("I1").Borders(10).LineStyle = 1
$ws.Range("I1").Borders(8).LineStyle = -4142
$ws.Range("I1").Borders(9).LineStyle = -4142

# Match the recorded selection/active cell.
[void]$ws.Range("I1").Select()
